$wb = $excel.ActiveWorkbook

# 1. Rename the second sheet
$ws2 = $wb.Worksheets.Item("Include from Ferlab.bio CodeS")
$ws2.Name = "Include #0"

# 2. Update Metadata sheet
$ws1 = $wb.Worksheets.Item("Metadata")

# a. Date value (row 8, col B)
$ws1.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# b. Contact value (row 10, col B)
$ws1.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# c. Insert new row after row 10 (Contact) for Jurisdiction
$ws1.Rows.Item(11).Insert()
$ws1.Range("A11").Value = "Jurisdiction"
# Use a leading apostrophe so the cell stores a real (typed) empty string,
# rather than being left as a truly blank/uninitialized cell.
$ws1.Range("B11").Value = "'"

# Copy style from row 10 to new row 11 so formatting matches other data rows
# (also clears the quote-prefix flag picked up from the "'" entry above).
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
